# Apply the cryptos-list refresh described by the commit:
#   "Updated cryptos list on Tue Aug  8 08:18:54 UTC 2023 with GitHub Actions"
#
# For each affected row we update Price (D) / Volume 1h (E), and for the
# four rows whose rank order swapped we also rewrite Coin (B) and Link (C).
# Numeric-looking Price strings are forced back to Text (NumberFormat "@")
# before the assignment so Excel does not silently reinterpret them as
# numbers (e.g. "1.000" -> 1, "0.9981" -> 0.9981 as a float) - the sheet
# stores every Price/Volume cell as plain text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.159.08"
$ws.Range("E2").Value = "  +0.44%  "

# Row 3
$ws.Range("D3").Value = "1.828.79"
$ws.Range("E3").Value = "  -0.25%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9981"
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.65"
$ws.Range("E5").Value = "  +0.27%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6205"
$ws.Range("E6").Value = "  +1.27%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07350"
$ws.Range("E8").Value = "  -1.58%  "

# Row 9
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2912"
$ws.Range("E9").Value = "  -0.28%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.23"
$ws.Range("E10").Value = "  +0.58%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07635"
$ws.Range("E11").Value = "  -0.58%  "

# Row 12
$ws.Range("D12").Value = "1.839.21"
$ws.Range("E12").Value = "  +0.25%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.964"
$ws.Range("E13").Value = "  -0.68%  "

# Row 14
$ws.Range("E14").Value = "  -0.48%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.35"
$ws.Range("E15").Value = "  -0.17%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009005"
$ws.Range("E16").Value = "  -1.58%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.838"
$ws.Range("E17").Value = "  -1.17%  "

# Row 18
$ws.Range("D18").Value = "29.156.66"
$ws.Range("E18").Value = "  +0.51%  "

# Row 19
$ws.Range("D19").Value = "2.085.68"
$ws.Range("E19").Value = "  -0.44%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "236.11"
$ws.Range("E20").Value = "  +2.43%  "

# Row 21
$ws.Range("E21").Value = "  -1.48%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9996"
$ws.Range("E22").Value = "  -0.08%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.354"
$ws.Range("E23").Value = "  +2.03%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  +0.02%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.16"
$ws.Range("E25").Value = "  -0.48%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1391"
$ws.Range("E26").Value = "  +0.01%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.528"
$ws.Range("E27").Value = "  +0.49%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.61"
$ws.Range("E28").Value = "  -0.82%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.489"
$ws.Range("E29").Value = "  -0.28%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05870"
$ws.Range("E30").Value = "  +6.45%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.079"
$ws.Range("E31").Value = "  -0.95%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.080"
$ws.Range("E32").Value = "  -1.62%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.211"
$ws.Range("E33").Value = "  +0.94%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.858"
$ws.Range("E34").Value = "  +1.13%  "

# Row 35
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7250"
$ws.Range("E35").Value = "  -2.03%  "

# Row 36
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.137"
$ws.Range("E36").Value = "  -0.31%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.608"
$ws.Range("E37").Value = "  -1.79%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.845"
$ws.Range("E38").Value = "  +2.77%  "

# Row 39
$ws.Range("D39").Value = "1.224.01"
$ws.Range("E39").Value = "  +1.26%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01754"
$ws.Range("E40").Value = "  -1.29%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.227"
$ws.Range("E41").Value = "  -3.74%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9073"
$ws.Range("E42").Value = "  +2.21%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  +0.13%  "

# Row 44
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.80"
$ws.Range("E44").Value = "  -0.03%  "

# Row 45
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.989.48"
$ws.Range("E45").Value = "  +0.00%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.53"
$ws.Range("E46").Value = "  +0.19%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5039"
$ws.Range("E47").Value = "  -0.87%  "

# Row 48
$ws.Range("E48").Value = "  -4.19%  "

# Row 49
$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4027"
$ws.Range("E49").Value = "  -0.87%  "

# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.090"
$ws.Range("E50").Value = "  -0.40%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1128"
$ws.Range("E51").Value = "  +2.62%  "
